# Auto-generated update of crypto Price (D) / Volume(1h) (E) columns for rows 2-51,
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    [PSCustomObject]@{ Row = 2; D = "44.292.29"; E = "  +2.59%  " }
    [PSCustomObject]@{ Row = 3; D = "2.431.86"; E = "  +2.12%  " }
    [PSCustomObject]@{ Row = 4; D = "1.00"; E = "  -0.04%  " }
    [PSCustomObject]@{ Row = 5; D = "307.92"; E = "  +1.55%  " }
    [PSCustomObject]@{ Row = 6; D = "99.96"; E = "  +3.12%  " }
    [PSCustomObject]@{ Row = 7; D = $null; E = "  +0.94%  " }
    [PSCustomObject]@{ Row = 8; D = $null; E = "  -0.03%  " }
    [PSCustomObject]@{ Row = 9; D = "0.500"; E = "  -0.22%  " }
    [PSCustomObject]@{ Row = 10; D = "35.45"; E = "  +3.87%  " }
    [PSCustomObject]@{ Row = 11; D = "0.0801"; E = "  +1.47%  " }
    [PSCustomObject]@{ Row = 12; D = $null; E = "  +2.66%  " }
    [PSCustomObject]@{ Row = 13; D = "18.77"; E = "  +1.64%  " }
    [PSCustomObject]@{ Row = 14; D = "6.94"; E = "  +2.40%  " }
    [PSCustomObject]@{ Row = 15; D = "2.808.97"; E = "  +2.09%  " }
    [PSCustomObject]@{ Row = 16; D = "2.479.69"; E = "  +5.01%  " }
    [PSCustomObject]@{ Row = 17; D = "0.832"; E = "  +3.03%  " }
    [PSCustomObject]@{ Row = 18; D = "44.264.22"; E = "  +2.57%  " }
    [PSCustomObject]@{ Row = 19; D = "12.37"; E = "  +1.16%  " }
    [PSCustomObject]@{ Row = 20; D = "6.46"; E = "  +1.58%  " }
    [PSCustomObject]@{ Row = 21; D = "0.0₃0908"; E = "  +1.99%  " }
    [PSCustomObject]@{ Row = 22; D = "68.69"; E = "  +0.23%  " }
    [PSCustomObject]@{ Row = 23; D = "240.70"; E = "  +2.21%  " }
    [PSCustomObject]@{ Row = 24; D = "2.30"; E = "  +3.80%  " }
    [PSCustomObject]@{ Row = 25; D = "2.48"; E = "  +1.80%  " }
    [PSCustomObject]@{ Row = 26; D = $null; E = "  +0.00%  " }
    [PSCustomObject]@{ Row = 27; D = "25.35"; E = "  +2.18%  " }
    [PSCustomObject]@{ Row = 28; D = "2.35"; E = "  -0.99%  " }
    [PSCustomObject]@{ Row = 29; D = "9.52"; E = "  +4.32%  " }
    [PSCustomObject]@{ Row = 30; D = "32.98"; E = "  +4.43%  " }
    [PSCustomObject]@{ Row = 31; D = $null; E = "  +16.85%  " }
    [PSCustomObject]@{ Row = 32; D = $null; E = "  +8.83%  " }
    [PSCustomObject]@{ Row = 33; D = "5.18"; E = "  +1.72%  " }
    [PSCustomObject]@{ Row = 34; D = $null; E = "  -0.08%  " }
    [PSCustomObject]@{ Row = 35; D = "0.0766"; E = "  +4.02%  " }
    [PSCustomObject]@{ Row = 36; D = "1.92"; E = "  +3.75%  " }
    [PSCustomObject]@{ Row = 37; D = "4.59"; E = "  +6.54%  " }
    [PSCustomObject]@{ Row = 38; D = "130.69"; E = "  +21.36%  " }
    [PSCustomObject]@{ Row = 39; D = "2.93"; E = "  +4.77%  " }
    [PSCustomObject]@{ Row = 40; D = $null; E = "  -0.69%  " }
    [PSCustomObject]@{ Row = 41; D = $null; E = "  +0.27%  " }
    [PSCustomObject]@{ Row = 42; D = "21.07"; E = "  -5.73%  " }
    [PSCustomObject]@{ Row = 43; D = "0.0288"; E = "  +2.76%  " }
    [PSCustomObject]@{ Row = 44; D = "1.958.75"; E = "  +0.17%  " }
    [PSCustomObject]@{ Row = 45; D = $null; E = "  +2.00%  " }
    [PSCustomObject]@{ Row = 46; D = "2.89"; E = "  +4.91%  " }
    [PSCustomObject]@{ Row = 47; D = "9.36"; E = "  +1.30%  " }
    [PSCustomObject]@{ Row = 48; D = $null; E = "  +9.48%  " }
    [PSCustomObject]@{ Row = 49; D = "2.661.84"; E = "  +1.84%  " }
    [PSCustomObject]@{ Row = 50; D = "53.59"; E = "  +1.45%  " }
    [PSCustomObject]@{ Row = 51; D = "73.76"; E = "  +2.36%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Price column holds formatted text (e.g. thousands-dot grouping, trailing
        # zeros); force Text format first so Excel does not reinterpret it as a number.
        $ws.Cells.Item($u.Row, 4).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
